$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("A5").Value = "CNN+LSTM+DNN "
$ws.Range("B5").Value = 0.784
$ws.Range("C5").Value = 0.765
$ws.Range("D5").Value = 0.701
$ws.Range("A6").Value = "BERT"
$ws.Range("B6").Value = 0.804
$ws.Range("C6").Value = 0.797
$ws.Range("D6").Value = 0.762
